$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.506.83'
$ws.Range("E2").Value = '  +0.01%  '

$ws.Range("D3").Value = '2.285.40'
$ws.Range("E3").Value = '  -0.40%  '

$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").Value = '311.24'
$ws.Range("E5").Value = '  -3.69%  '

$ws.Range("D6").Value = '103.32'
$ws.Range("E6").Value = '  -1.31%  '

$ws.Range("D7").Value = '0.623'
$ws.Range("E7").Value = '  -1.07%  '

$ws.Range("E8").Value = '  +0.13%  '

$ws.Range("E9").Value = '  -1.17%  '

$ws.Range("D10").Value = '38.71'
$ws.Range("E10").Value = '  -3.86%  '

$ws.Range("D11").Value = '0.0901'
$ws.Range("E11").Value = '  -0.86%  '

$ws.Range("D12").Value = '8.17'
$ws.Range("E12").Value = '  -4.16%  '

$ws.Range("D13").Value = '0.108'
$ws.Range("E13").Value = '  +0.53%  '

$ws.Range("D14").Value = '0.970'
$ws.Range("E14").Value = '  -0.17%  '

$ws.Range("D15").Value = '15.23'
$ws.Range("E15").Value = '  -0.06%  '

$ws.Range("D16").Value = '2.631.50'
$ws.Range("E16").Value = '  -0.48%  '

$ws.Range("D17").Value = '2.285.38'
$ws.Range("E17").Value = '  -0.60%  '

$ws.Range("D18").Value = '42.710.83'
$ws.Range("E18").Value = '  +0.17%  '

$ws.Range("D19").Value = '7.31'
$ws.Range("E19").Value = '  -1.94%  '

$ws.Range("B20").Value = 'InternetComputer(DFINITY)'
$ws.Range("C20").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D20").Value = '13.58'
$ws.Range("E20").Value = '  +3.46%  '

$ws.Range("B21").Value = 'ShibaInu'
$ws.Range("C21").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D21").Value = '0.0000104'
$ws.Range("E21").Value = '  -1.89%  '

$ws.Range("D22").Value = '73.19'
$ws.Range("E22").Value = '  -0.32%  '

$ws.Range("D23").Value = '268.59'
$ws.Range("E23").Value = '  -0.78%  '

$ws.Range("E24").Value = '  -5.10%  '

$ws.Range("D25").Value = '2.17'
$ws.Range("E25").Value = '  -2.81%  '

$ws.Range("D26").Value = '1.01'
$ws.Range("E26").Value = '  +0.11%  '

$ws.Range("D27").Value = '10.74'
$ws.Range("E27").Value = '  -1.72%  '

$ws.Range("D28").Value = '7.09'
$ws.Range("E28").Value = '  +15.83%  '

$ws.Range("D29").Value = '2.30'
$ws.Range("E29").Value = '  -1.68%  '

$ws.Range("D30").Value = '22.34'
$ws.Range("E30").Value = '  -1.16%  '

$ws.Range("D31").Value = '35.65'
$ws.Range("E31").Value = '  -7.67%  '

$ws.Range("D32").Value = '164.59'
$ws.Range("E32").Value = '  -0.43%  '

$ws.Range("D33").Value = '0.0848'
$ws.Range("E33").Value = '  -4.00%  '

$ws.Range("D34").Value = '0.130'
$ws.Range("E34").Value = '  -2.42%  '

$ws.Range("E35").Value = '  +0.64%  '

$ws.Range("D36").Value = '0.112'
$ws.Range("E36").Value = '  -2.79%  '

$ws.Range("E37").Value = '  -3.06%  '

$ws.Range("D38").Value = '0.0346'
$ws.Range("E38").Value = '  -2.63%  '

$ws.Range("E39").Value = '  +0.67%  '

$ws.Range("D40").Value = '3.60'
$ws.Range("E40").Value = '  -4.14%  '

$ws.Range("D41").Value = '111.62'
$ws.Range("E41").Value = '  +16.17%  '

$ws.Range("D42").Value = '1.56'
$ws.Range("E42").Value = '  +1.33%  '

$ws.Range("D43").Value = '70.52'
$ws.Range("E43").Value = '  +0.41%  '

$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").Value = '0.225'
$ws.Range("E44").Value = '  -0.05%  '

$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").Value = '0.999'
$ws.Range("E45").Value = '  -0.45%  '

$ws.Range("D46").Value = '12.06'
$ws.Range("E46").Value = '  -2.20%  '

$ws.Range("D47").Value = '1.723.86'
$ws.Range("E47").Value = '  +8.75%  '

$ws.Range("D48").Value = '110.38'
$ws.Range("E48").Value = '  -3.21%  '

$ws.Range("D49").Value = '77.36'
$ws.Range("E49").Value = '  -5.10%  '

$ws.Range("D50").Value = '8.66'
$ws.Range("E50").Value = '  -2.86%  '

$ws.Range("D51").Value = '5.15'
$ws.Range("E51").Value = '  -2.59%  '
